$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header labels in A1/B1 (lat/long -> long/lat)
$ws.Range("A1").Value = "long"
$ws.Range("B1").Value = "lat"

# Row 2 - swap A/B values and update C/D
$ws.Range("A2").Value = -35.19630277777777
$ws.Range("B2").Value = -8.041286111111111
$ws.Range("C2").Value = 66.26624612895462
$ws.Range("D2").Value = 173.2656707633173

# Row 3
$ws.Range("C3").Value = -25.55414410611931
$ws.Range("D3").Value = 142.6785961197031

# Row 4
$ws.Range("C4").Value = -163.6055372995401
$ws.Range("D4").Value = 212.6335549214862

# Row 5
$ws.Range("C5").Value = 9.468649267773118
$ws.Range("D5").Value = 234.5252131391234

# Row 6
$ws.Range("C6").Value = 129.0682053998726
$ws.Range("D6").Value = 43.38276794851311

# Row 7
$ws.Range("C7").Value = 118.0197971337859
$ws.Range("D7").Value = 102.6701866036969

# Row 8
$ws.Range("C8").Value = 27.99423961891065
$ws.Range("D8").Value = 62.62552101450152

# Row 9
$ws.Range("C9").Value = 57.82247016304152
$ws.Range("D9").Value = 42.68871649655923

# Row 10
$ws.Range("C10").Value = -141.7638576015277
$ws.Range("D10").Value = 115.7712185651046

# Row 11
$ws.Range("C11").Value = 87.00951302332027
$ws.Range("D11").Value = 200.2221681110879

# Row 12
$ws.Range("C12").Value = -115.2738503653425
$ws.Range("D12").Value = 34.53867755091041

# Row 13
$ws.Range("C13").Value = -121.6453102792624
$ws.Range("D13").Value = 398.1675225514543

# Row 14
$ws.Range("C14").Value = -11.31302706143295
$ws.Range("D14").Value = 174.8205722261661

# Row 15
$ws.Range("C15").Value = -52.62463068105861
$ws.Range("D15").Value = 203.4467576568445

# Row 16
$ws.Range("C16").Value = 136.6447570359038
$ws.Range("D16").Value = 264.4902170227521

# Row 17 - swap A/B values and update C/D
$ws.Range("A17").Value = -35.19621944444444
$ws.Range("B17").Value = -8.041044444444445
$ws.Range("C17").Value = -48.15074066962603
$ws.Range("D17").Value = 191.554550469853

# Row 18
$ws.Range("C18").Value = -8.786703908889159
$ws.Range("D18").Value = 116.2614160447615

# Row 19
$ws.Range("C19").Value = -15.296833632902
$ws.Range("D19").Value = 189.185887616363

# Row 20
$ws.Range("C20").Value = -33.66904861951783
$ws.Range("D20").Value = 148.0218177577202

# Row 21
$ws.Range("C21").Value = -15.74843470576612
$ws.Range("D21").Value = 196.303535000873

# Row 22
$ws.Range("C22").Value = -48.05156290250891
$ws.Range("D22").Value = 138.3335131970775

# Row 23
$ws.Range("C23").Value = 139.9589635289357
$ws.Range("D23").Value = 260.8294987402039

# Row 24
$ws.Range("C24").Value = 172.7179099406849
$ws.Range("D24").Value = 132.8652132539594

# Row 25
$ws.Range("C25").Value = 102.6511282407744
$ws.Range("D25").Value = 92.5778163111207

# Row 26
$ws.Range("C26").Value = 136.2801747462221
$ws.Range("D26").Value = 178.1047124542434

# Row 27
$ws.Range("C27").Value = 136.7715862000828
$ws.Range("D27").Value = 207.44112556738

# Row 28
$ws.Range("C28").Value = -148.4456885486338
$ws.Range("D28").Value = 111.7596901200424

# Row 29
$ws.Range("C29").Value = -71.16438397036443
$ws.Range("D29").Value = 32.3526590805841

# Row 30
$ws.Range("C30").Value = 102.56878929778
$ws.Range("D30").Value = 36.70518828241583

# Row 31
$ws.Range("C31").Value = 66.26624612895462
$ws.Range("D31").Value = 173.2656707633173
